$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The oldest reporting period (column D, "6 ماهه منتهی به 1399/06" /
# "1400-09-01 (4)") is dropped, and every later period shifts one column
# to the left. Deleting column D performs that shift for all rows
# (values, styles, shared-string references, column widths) in one go.
$ws.Columns("D").Delete()

# After the shift, column M (13) is empty. Re-create it by duplicating
# column L (12) - both its value and its formatting - for every data row,
# then overwrite the cells that actually hold new figures for the newest
# period (12 ماهه منتهی به 1401/12).
for ($r = 8; $r -le 27; $r++) {
    $ws.Cells.Item($r, 12).Copy()
    $ws.Cells.Item($r, 13).PasteSpecial(-4122)
    $ws.Cells.Item($r, 13).Value2 = $ws.Cells.Item($r, 12).Value2
}

# New column headers for the newly-added period / publish date.
$ws.Cells.Item(8, 13).Value2 = "12 ماهه منتهی به 1401/12"
$ws.Cells.Item(9, 13).Value2 = "1402-02-25"

# The most recent previously-existing publish-date label is revised
# in place (it shifted from column J to column I after the delete).
$ws.Cells.Item(9, 9).Value2 = "1402-02-25 (8)"

# Updated cumulative figures for the newest period (column M).
$ws.Cells.Item(11, 13).Value2 = 77250
$ws.Cells.Item(12, 13).Value2 = -42971
$ws.Cells.Item(13, 13).Value2 = 34279
$ws.Cells.Item(14, 13).Value2 = -13504
$ws.Cells.Item(17, 13).Value2 = 20775
$ws.Cells.Item(18, 13).Value2 = -2943
$ws.Cells.Item(19, 13).Value2 = 287
$ws.Cells.Item(20, 13).Value2 = 18119
$ws.Cells.Item(21, 13).Value2 = -725
$ws.Cells.Item(22, 13).Value2 = 17394
$ws.Cells.Item(24, 13).Value2 = 17394
$ws.Cells.Item(26, 13).Value2 = 11335

Write-Host "Edit complete"
